{"js": "const body = context.document.body;\n\n// Locate the existing list-item paragraph that references\n// \"Hamburger_icon.svg\" \u2014 the new Brawlhalla entry goes right after it\n// (and before the trailing blank list item already in the document).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Hamburger_icon.svg\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not find the Hamburger_icon.svg reference paragraph\");\n}\n\n// Insert a brand-new paragraph right after it; Word copies the paragraph\n// formatting (ListParagraph style, numbering, tab stops) automatically.\nconst newParagraph = targetParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// Type the URL, turn that exact range into a hyperlink, then append the\n// trailing space as its own run (matching the source run split).\nconst url = \"https://www.brawlhalla.com/about/\";\nconst urlRange = newParagraph.insertText(url, \"Start\");\nawait context.sync();\n\nurlRange.hyperlink = url;\nawait context.sync();\n\nnewParagraph.insertText(\" \", \"End\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the existing list-item paragraph that contains the\n# \"Hamburger_icon.svg\" reference hyperlink - the new \"Brawlhalla\" entry\n# is inserted immediately after it (and before the trailing blank list item).\n$targetPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*Hamburger_icon.svg*\") {\n    $targetPara = $p\n  }\n}\n\n# Insert a brand-new paragraph right after it; Word copies the paragraph\n# formatting (ListParagraph style, numbering, tab stops) from $targetPara.\n$targetPara.Range.InsertParagraphAfter()\n\n# Re-fetch the freshly minted paragraph (index = old index + 1).\n$newPara = $d.Paragraphs.Item($targetPara.Index + 1)\n$newRange = $newPara.Range\n\n# Type the URL text followed by a trailing space, matching the source run.\n$url = \"https://www.brawlhalla.com/about/\"\n$newRange.InsertAfter($url + \" \")\n\n# Turn just the URL portion (not the trailing space) into a real hyperlink.\n$linkRange = $d.Range($newRange.Start, $newRange.Start + $url.Length)\n$d.Hyperlinks.Add($linkRange, $url) | Out-Null\n"}
